# Spelling/accent corrections on the "Inspeccion de código" checklist sheet
# (commit: "check list ortogrfia corregida").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inspeccion de código")

# Observaciones text correction (ususario -> usuario, categoria -> categoría)
$ws.Range("G8").Value = "El descuento es ingresado de manera manual por el usuario en ves de ser en base  a la categoría del producto"

# Inspector name: "Matias" -> "Matías"
$ws.Range("C5").Value = "Inspector: Matías Fuentealba, Vicente Zurita"

# Header label "Nro" -> "Nro."
$ws.Range("B7").Value = "Nro."

# Method name accent: "Categoria" -> "Categoría"
$ws.Range("C8").Value = "setDescuento Categoría()"

# Remaining Observaciones text corrections (metodo -> método, informacion -> información, estatico -> estático)
$ws.Range("G9").Value = "Este método no se pide en modelo de clases"
$ws.Range("G10").Value = "Realiza correctamente el calculo, pero no existe un método que imprima la información"
$ws.Range("G11").Value = "Este método solo muestra el atributo ''tipo`""
$ws.Range("G13").Value = "Este método utiliza un valor de IVA estático en vez de los valores de la interfaz y monto "
$ws.Range("G14").Value = "Esta clase contiene métodos que no se usan ni se piden"

# Move the active selection from C25 to C20, matching the saved view state.
$ws.Range("C20").Select()
